# Weekly Albahaca (Mercado Mayorista Lo Valledor de Santiago) price-list
# update: a new weekly observation is inserted as row 670, pushing all
# subsequent observations down by one row (old row 670 -> 671, ...,
# old row 771 -> 772). Sheet used range grows from A1:R771 to A1:R772.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 670; Excel shifts rows 670:771 down to 671:772
# and extends the sheet dimension automatically.
$ws.Rows(670).Insert()

# Populate the newly inserted row 670 with the new observation.
$ws.Range("A670").Value = 6
$ws.Range("B670").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C670").Value = "Metropolitana"
$ws.Range("D670").Value = 45127
$ws.Range("E670").Value = 13
$ws.Range("F670").Value = 100112052
$ws.Range("G670").Value = "Albahaca"
$ws.Range("H670").Value = "Sin especificar"
$ws.Range("I670").Value = "Primera"
$ws.Range("J670").Value = 80
$ws.Range("K670").Value = 4000
$ws.Range("L670").Value = 4500
$ws.Range("M670").Value = 4219
$ws.Range("N670").Value = "`$/paquete"
$ws.Range("O670").Value = "Región de Arica y Parinacota"
$ws.Range("P670").Value = 4219
$ws.Range("Q670").Value = 1
$ws.Range("R670").Value = "Hortaliza"
